# Update the "chi_nhanh" (branch) column on the KhachHang sheet.
# The demo branch names (city-based) are being replaced with actual
# district names used by the system (Tra Vinh province districts).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KhachHang")

$ws.Range("C2").Value = "Càng Long"
$ws.Range("C3").Value = "Càng Long"
$ws.Range("C4").Value = "Cầu Kè"
$ws.Range("C5").Value = "Cầu Kè"
$ws.Range("C6").Value = "Tiểu Cần"
$ws.Range("C7").Value = "Châu Thành"
$ws.Range("C8").Value = "Trà Cú"
$ws.Range("C9").Value = "Trà Cú"
